# Tuntikirjanpito.xlsx update
# - add a new work-log entry (row 84): 18.1.2022, 1h,
#   "punta lisätty valuuttoihin, uusien komponenttien siistimistä, analyzer
#    työkalun css yhdenmukaiseksi, coin logot"
# - the three summary rows (tunnit yht. / target / suoritettu(%)) move down
#   to make room and their formulas are recomputed to cover the new row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) new log entry in row 84 --------------------------------------------------
$ws.Range("A15").Copy($ws.Range("A84"))
$ws.Range("A84").Value = 44579

$ws.Range("B2").Copy($ws.Range("B84"))
$ws.Range("B84").Value = 1

$ws.Range("C84").Value = "punta lisätty valuuttoihin, uusien komponenttien siistimistä, analyzer työkalun css yhdenmukaiseksi, coin logot"

# --- 2) write the summary block to rows 99-101; they will land on 96-98 once the
#        old summary rows (86-88) are deleted a few lines down ------------------
$ws.Range("A86").Copy($ws.Range("A99"))
$ws.Range("B99").Formula = "=SUM(B2:B84)"
$ws.Rows.Item(99).RowHeight = 14.25

$ws.Range("A87").Copy($ws.Range("A100"))
$ws.Range("B87").Copy($ws.Range("B100"))

$ws.Range("A88").Copy($ws.Range("A101"))
$ws.Range("B101").Formula = "=B99/B100*100"

# --- 3) drop the old summary rows; this shifts 99-101 up onto 96-98 -------------
$ws.Range("86:88").Delete()

# --- 4) view state: active cell + scroll position -------------------------------
$null = $ws.Range("C85").Select()
